$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 14 : 07/08/2021 Sat
# ---------------------------------------------------------------------------
$ws.Range("A14").Value = "07/08/2021 Sat"
$ws.Range("B14").Value = "NA"
$ws.Range("C14").Value = "NA"
$ws.Range("D14").Value = "NA"
$ws.Range("G14").Value = "solving many problems"

# ---------------------------------------------------------------------------
# Row 15 : 8/8/2021 Sunday  (A15 uses the same date number format as A9/A13)
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A15").Value = "8/8/2021 Sunday"
$ws.Range("B15").Value = "2 Backtracking problem "
$ws.Range("C15").Value = "nA"
$ws.Range("D15").Value = "Na"
$ws.Range("E15").Value = "NA"
$ws.Range("G15").Value = "Solvong appti problem"

# ---------------------------------------------------------------------------
# Row 16 : 09/08/2021 Monday
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = "09/08/2021 Monday"
$ws.Range("B16").Value = "2 problem on back tracking "
$ws.Range("C16").Value = "give contest of DIV 2"
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "``"
$ws.Range("G16").Value = "Solving ratio problem"

# ---------------------------------------------------------------------------
# View state: scroll / active cell moves to G16
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$null = $ws.Range("G16").Select()
